$d = $word.ActiveDocument

# The last paragraph contains an inline picture (run 1) followed by an
# empty run (run 2, rPr rtl=0). The edit removes the picture and puts
# the text "Não há restrições." into that trailing empty run.

$shape = $d.InlineShapes.Item(1)
$shape.Delete()

$lastParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastParagraph.Range.Text = "Não há restrições."
